# Updated symbol list on Sat Jan 14 22:43:18 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (column D) and "Volume(1h)" (column E) values for the
# crypto symbols whose quotes changed. All of these cells store plain text
# (e.g. "303.79", "3.51%") rather than real numbers/percentages, so a leading
# apostrophe is used to force Excel to keep them as literal text instead of
# auto-converting them to numeric/percentage values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.79"
$ws.Range("E2").Value = "'3.51%"

$ws.Range("D3").Value = "'32.05"
$ws.Range("E3").Value = "'8.64%"

$ws.Range("D4").Value = "'5.236"
$ws.Range("E4").Value = "'0.08%"

$ws.Range("D5").Value = "'0.07546"
$ws.Range("E5").Value = "'5.30%"

$ws.Range("D6").Value = "'7.953"
$ws.Range("E6").Value = "'5.34%"

$ws.Range("D7").Value = "'3.847"
$ws.Range("E7").Value = "'7.07%"

$ws.Range("D8").Value = "'1.528"
$ws.Range("E8").Value = "'8.72%"

$ws.Range("D9").Value = "'0.9287"
$ws.Range("E9").Value = "'2.17%"

$ws.Range("D10").Value = "'0.1695"
$ws.Range("E10").Value = "'3.97%"

$ws.Range("D11").Value = "'0.07903"
$ws.Range("E11").Value = "'2.01%"

$ws.Range("D12").Value = "'0.08011"
$ws.Range("E12").Value = "'3.32%"

$ws.Range("D13").Value = "'0.03036"
$ws.Range("E13").Value = "'3.98%"

$ws.Range("D14").Value = "'0.09903"
$ws.Range("E14").Value = "'10.00%"

$ws.Range("D15").Value = "'0.001494"
$ws.Range("E15").Value = "'-7.31%"

$ws.Range("D16").Value = "'0.04592"
$ws.Range("E16").Value = "'1.35%"

$ws.Range("D17").Value = "'0.006518"
$ws.Range("E17").Value = "'7.23%"

$ws.Range("D18").Value = "'3.446"
$ws.Range("E18").Value = "'-1.16%"

$ws.Range("D19").Value = "'2.228"
$ws.Range("E19").Value = "'-0.24%"

$ws.Range("E20").Value = "'1.39%"

$ws.Range("D21").Value = "'0.1334"
$ws.Range("E21").Value = "'-2.39%"

$ws.Range("D22").Value = "'4.457"
$ws.Range("E22").Value = "'10.10%"

$ws.Range("D23").Value = "'0.1618"
$ws.Range("E23").Value = "'1.71%"

$ws.Range("D24").Value = "'0.001215"
$ws.Range("E24").Value = "'0.75%"

$ws.Range("D25").Value = "'0.004471"
$ws.Range("E25").Value = "'4.89%"

$ws.Range("D26").Value = "'0.0001396"
$ws.Range("E26").Value = "'20.13%"

$ws.Range("D27").Value = "'0.0001782"
$ws.Range("E27").Value = "'5.85%"

$ws.Range("D39").Value = "'0.01693"
$ws.Range("E39").Value = "'2,486.60%"

$ws.Range("D40").Value = "'0.04497"
$ws.Range("E40").Value = "'1.41%"

$ws.Range("D41").Value = "'0.006985"
$ws.Range("E41").Value = "'-0.21%"

$ws.Range("E42").Value = "'5.68%"

$ws.Range("D43").Value = "'0.002075"
$ws.Range("E43").Value = "'-5.50%"

$ws.Range("D44").Value = "'0.01371"
$ws.Range("E44").Value = "'3.28%"

$ws.Range("D45").Value = "'0.00006164"
$ws.Range("E45").Value = "'5.55%"

$ws.Range("D46").Value = "'0.7191"
$ws.Range("E46").Value = "'-62.73%"

$ws.Range("D47").Value = "'0.01297"
$ws.Range("E47").Value = "'0.19%"
